$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '98.509.35'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.455.11'
$ws.Range('E3').Value = '  +4.15%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '260.60'
$ws.Range('E5').Value = '  +1.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '670.57'
$ws.Range('E6').Value = '  +7.43%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.56'
$ws.Range('E7').Value = '  +8.62%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.456'
$ws.Range('E8').Value = '  +12.85%  '
$ws.Range('E9').Value = '  +21.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.999'
$ws.Range('E10').Value = '  +0.00%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '3.448.43'
$ws.Range('E11').Value = '  +4.01%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.218'
$ws.Range('E12').Value = '  +9.27%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '43.21'
$ws.Range('E13').Value = '  +9.75%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000271'
$ws.Range('E14').Value = '  +8.32%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.24'
$ws.Range('E15').Value = '  +13.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '98.163.10'
$ws.Range('E16').Value = '  -0.40%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.103.80'
$ws.Range('E17').Value = '  +4.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.84'
$ws.Range('E18').Value = '  +39.61%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.456.36'
$ws.Range('E19').Value = '  +4.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.65'
$ws.Range('E20').Value = '  +14.98%  '
$ws.Range('B21').Value = 'Stellar'
$ws.Range('C21').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.511'
$ws.Range('E21').Value = '  +74.88%  '
$ws.Range('B22').Value = 'SuiNetwork'
$ws.Range('C22').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.58'
$ws.Range('E22').Value = '  +2.43%  '
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '524.43'
$ws.Range('E23').Value = '  +8.22%  '
$ws.Range('B24').Value = 'Uniswap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.92'
$ws.Range('E24').Value = '  +15.58%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000215'
$ws.Range('E25').Value = '  +4.86%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.42'
$ws.Range('E26').Value = '  +14.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '103.16'
$ws.Range('E27').Value = '  +16.52%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '13.04'
$ws.Range('E28').Value = '  +8.53%  '
$ws.Range('B29').Value = 'WrappedeETH'
$ws.Range('C29').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.653.16'
$ws.Range('E29').Value = '  +4.62%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.156'
$ws.Range('E30').Value = '  +19.61%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '11.73'
$ws.Range('E31').Value = '  +14.40%  '
$ws.Range('B32').Value = 'Cronos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.197'
$ws.Range('E32').Value = '  +4.79%  '
$ws.Range('B33').Value = 'Dai'
$ws.Range('C33').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.998'
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('B34').Value = 'PolygonEcosystemToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.603'
$ws.Range('E34').Value = '  +30.14%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.01'
$ws.Range('E35').Value = '  +0.47%  '
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '30.40'
$ws.Range('E36').Value = '  +8.72%  '
$ws.Range('B37').Value = 'PancakeSwap'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.25'
$ws.Range('E37').Value = '  +15.21%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.07'
$ws.Range('E38').Value = '  +11.32%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.161'
$ws.Range('E39').Value = '  +8.43%  '
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '534.03'
$ws.Range('E40').Value = '  +8.56%  '
$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.42'
$ws.Range('E41').Value = '  +14.16%  '
$ws.Range('B42').Value = 'WhiteBITCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '24.76'
$ws.Range('E42').Value = '  -0.25%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0449'
$ws.Range('E43').Value = '  +35.78%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.864'
$ws.Range('E44').Value = '  +8.45%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.75'
$ws.Range('E45').Value = '  +3.57%  '
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.45'
$ws.Range('E46').Value = '  +9.40%  '
$ws.Range('B47').Value = 'Cosmos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.59'
$ws.Range('E47').Value = '  +15.18%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.39'
$ws.Range('E48').Value = '  +14.67%  '
$ws.Range('B49').Value = 'ImmutableX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.60'
$ws.Range('E49').Value = '  +17.47%  '
$ws.Range('B50').Value = 'USDe'
$ws.Range('C50').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.00'
$ws.Range('E50').Value = '  +0.04%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.11'
$ws.Range('E51').Value = '  +8.99%  '
